$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 150.0354306666667
$ws.Range("H2").Value = 450.106292
$ws.Range("I2").Value = 0.4152507364956075
$ws.Range("J2").Value = 0.4152507364956075
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 1.854571666666667
$ws.Range("N2").Value = 5.563715
$ws.Range("O2").Value = 0.01651371646154392
$ws.Range("P2").Value = 0.01651371646154392
$ws.Range("Q2").Value = 278.2514587105311
$ws.Range("R2").Value = 2504.26312839478
$ws.Range("S2").Value = 0.00685733292293575
$ws.Range("T2").Value = 0.00685733292293575

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 150.0354306666667
$ws.Range("H3").Value = 450.106292
$ws.Range("I3").Value = 0.4152507364956075
$ws.Range("J3").Value = 0.4152507364956075
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 92.91372433333333
$ws.Range("N3").Value = 278.741173
$ws.Range("O3").Value = 0.8273343794712995
$ws.Range("P3").Value = 0.8273343794712996
$ws.Range("Q3").Value = 13940.35064519561
$ws.Range("R3").Value = 125463.1558067605
$ws.Range("S3").Value = 0.3435512104035935
$ws.Range("T3").Value = 0.3435512104035935

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 150.0354306666667
$ws.Range("H4").Value = 450.106292
$ws.Range("I4").Value = 0.4152507364956075
$ws.Range("J4").Value = 0.4152507364956075
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 0.162136
$ws.Range("N4").Value = 0.4864080000000001
$ws.Range("O4").Value = 0.001443712303133186
$ws.Range("P4").Value = 0.001443712303133187
$ws.Range("Q4").Value = 24.32614458657067
$ws.Range("R4").Value = 218.935301279136
$ws.Range("S4").Value = 0.0005995025971638254
$ws.Range("T4").Value = 0.0005995025971638254

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 150.0354306666667
$ws.Range("H5").Value = 450.106292
$ws.Range("I5").Value = 0.4152507364956075
$ws.Range("J5").Value = 0.4152507364956075
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 17.37449166666667
$ws.Range("N5").Value = 52.123475
$ws.Range("O5").Value = 0.1547081917640233
$ws.Range("P5").Value = 0.1547081917640233
$ws.Range("Q5").Value = 2606.789339822744
$ws.Range("R5").Value = 23461.1040584047
$ws.Range("S5").Value = 0.06424269057191435
$ws.Range("T5").Value = 0.06424269057191435

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 68.382243
$ws.Range("H6").Value = 205.146729
$ws.Range("I6").Value = 0.1892604742946246
$ws.Range("J6").Value = 0.1892604742946246
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 1.854571666666667
$ws.Range("N6").Value = 5.563715
$ws.Range("O6").Value = 0.01651371646154392
$ws.Range("P6").Value = 0.01651371646154392
$ws.Range("Q6").Value = 126.819770370915
$ws.Range("R6").Value = 1141.377933338235
$ws.Range("S6").Value = 0.003125393809878753
$ws.Range("T6").Value = 0.003125393809878752

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 68.382243
$ws.Range("H7").Value = 205.146729
$ws.Range("I7").Value = 0.1892604742946246
$ws.Range("J7").Value = 0.1892604742946246
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 92.91372433333333
$ws.Range("N7").Value = 278.741173
$ws.Range("O7").Value = 0.8273343794712995
$ws.Range("P7").Value = 0.8273343794712996
$ws.Range("Q7").Value = 6353.648875397013
$ws.Range("R7").Value = 57182.83987857312
$ws.Range("S7").Value = 0.1565816970589871
$ws.Range("T7").Value = 0.1565816970589871

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 68.382243
$ws.Range("H8").Value = 205.146729
$ws.Range("I8").Value = 0.1892604742946246
$ws.Range("J8").Value = 0.1892604742946246
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 0.162136
$ws.Range("N8").Value = 0.4864080000000001
$ws.Range("O8").Value = 0.001443712303133186
$ws.Range("P8").Value = 0.001443712303133187
$ws.Range("Q8").Value = 11.087223351048
$ws.Range("R8").Value = 99.78501015943201
$ws.Range("S8").Value = 0.0002732376752359718
$ws.Range("T8").Value = 0.0002732376752359717

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 68.382243
$ws.Range("H9").Value = 205.146729
$ws.Range("I9").Value = 0.1892604742946246
$ws.Range("J9").Value = 0.1892604742946246
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 17.37449166666667
$ws.Range("N9").Value = 52.123475
$ws.Range("O9").Value = 0.1547081917640233
$ws.Range("P9").Value = 0.1547081917640233
$ws.Range("Q9").Value = 1188.106711151475
$ws.Range("R9").Value = 10692.96040036327
$ws.Range("S9").Value = 0.02928014575052279
$ws.Range("T9").Value = 0.02928014575052279

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 104.737245
$ws.Range("H10").Value = 314.211735
$ws.Range("I10").Value = 0.2898796499701289
$ws.Range("J10").Value = 0.2898796499701289
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 1.854571666666667
$ws.Range("N10").Value = 5.563715
$ws.Range("O10").Value = 0.01651371646154392
$ws.Range("P10").Value = 0.01651371646154392
$ws.Range("Q10").Value = 194.242727021725
$ws.Range("R10").Value = 1748.184543195525
$ws.Range("S10").Value = 0.004786990347578308
$ws.Range("T10").Value = 0.004786990347578308

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 104.737245
$ws.Range("H11").Value = 314.211735
$ws.Range("I11").Value = 0.2898796499701289
$ws.Range("J11").Value = 0.2898796499701289
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 92.91372433333333
$ws.Range("N11").Value = 278.741173
$ws.Range("O11").Value = 0.8273343794712995
$ws.Range("P11").Value = 0.8273343794712996
$ws.Range("Q11").Value = 9731.527509362793
$ws.Range("R11").Value = 87583.74758426515
$ws.Range("S11").Value = 0.2398274003293941
$ws.Range("T11").Value = 0.2398274003293941

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 104.737245
$ws.Range("H12").Value = 314.211735
$ws.Range("I12").Value = 0.2898796499701289
$ws.Range("J12").Value = 0.2898796499701289
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 0.162136
$ws.Range("N12").Value = 0.4864080000000001
$ws.Range("O12").Value = 0.001443712303133186
$ws.Range("P12").Value = 0.001443712303133187
$ws.Range("Q12").Value = 16.98167795532
$ws.Range("R12").Value = 152.83510159788
$ws.Range("S12").Value = 0.0004185028170898167
$ws.Range("T12").Value = 0.0004185028170898167

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 104.737245
$ws.Range("H13").Value = 314.211735
$ws.Range("I13").Value = 0.2898796499701289
$ws.Range("J13").Value = 0.2898796499701289
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 17.37449166666667
$ws.Range("N13").Value = 52.123475
$ws.Range("O13").Value = 0.1547081917640233
$ws.Range("P13").Value = 0.1547081917640233
$ws.Range("Q13").Value = 1819.756390442125
$ws.Range("R13").Value = 16377.80751397912
$ws.Range("S13").Value = 0.04484675647606666
$ws.Range("T13").Value = 0.04484675647606665

$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 38.15794
$ws.Range("H14").Value = 114.47382
$ws.Range("I14").Value = 0.105609139239639
$ws.Range("J14").Value = 0.105609139239639
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 1.854571666666667
$ws.Range("N14").Value = 5.563715
$ws.Range("O14").Value = 0.01651371646154392
$ws.Range("P14").Value = 0.01651371646154392
$ws.Range("Q14").Value = 70.76663438236668
$ws.Range("R14").Value = 636.8997094413
$ws.Range("S14").Value = 0.00174399938115111
$ws.Range("T14").Value = 0.00174399938115111

$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 38.15794
$ws.Range("H15").Value = 114.47382
$ws.Range("I15").Value = 0.105609139239639
$ws.Range("J15").Value = 0.105609139239639
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 92.91372433333333
$ws.Range("N15").Value = 278.741173
$ws.Range("O15").Value = 0.8273343794712995
$ws.Range("P15").Value = 0.8273343794712996
$ws.Range("Q15").Value = 3545.396318287874
$ws.Range("R15").Value = 31908.56686459086
$ws.Range("S15").Value = 0.0873740716793248
$ws.Range("T15").Value = 0.08737407167932479

$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 38.15794
$ws.Range("H16").Value = 114.47382
$ws.Range("I16").Value = 0.105609139239639
$ws.Range("J16").Value = 0.105609139239639
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 0.162136
$ws.Range("N16").Value = 0.4864080000000001
$ws.Range("O16").Value = 0.001443712303133186
$ws.Range("P16").Value = 0.001443712303133187
$ws.Range("Q16").Value = 6.186775759840001
$ws.Range("R16").Value = 55.68098183856001
$ws.Range("S16").Value = 0.0001524692136435726
$ws.Range("T16").Value = 0.0001524692136435726

$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 38.15794
$ws.Range("H17").Value = 114.47382
$ws.Range("I17").Value = 0.105609139239639
$ws.Range("J17").Value = 0.105609139239639
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 17.37449166666667
$ws.Range("N17").Value = 52.123475
$ws.Range("O17").Value = 0.1547081917640233
$ws.Range("P17").Value = 0.1547081917640233
$ws.Range("Q17").Value = 662.9748105471667
$ws.Range("R17").Value = 5966.7732949245
$ws.Range("S17").Value = 0.01633859896551951
$ws.Range("T17").Value = 0.01633859896551951
